# Update stats for 2025-09 (row 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B22").Value = 6289
$ws.Range("D22").Value = 5835120
$ws.Range("E22").Value = 927.8295436476387
$ws.Range("F22").Value = 8.263040110173868
$ws.Range("H22").Value = 26.89567584965074
